$d = $word.ActiveDocument

# --- Change 1: Heading "AutoPilot " (remove proofErr spell markers, merge runs) ---
$d.Content.Find.Execute("AutoPilot Pre-Provisioning ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AutoPilot Pre-Provisioning ", 2)

# --- Change 2: Body paragraph "...AutoPilot White Glove process..." (remove proofErr markers) ---
$d.Content.Find.Execute("This document is intended to create clarification on the AutoPilot White Glove process for Surface devices", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This document is intended to create clarification on the AutoPilot White Glove process for Surface devices", 2)

# --- Change 3: Rewrite the tenant/personnel sentence in Step 1 ---
$d.Content.Find.Execute("tenant. Since all the students are using the same tenant, these configurations should be done in the customer tenant. A personel can do the configuration with the guidance from Microsoft FastTrack Center,", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "tenant. This configuration includes user information, device information and configuration information such as certificates, applications, or profiles so it should be organizations production tenant. A personnel can do the configuration with the guidance from Microsoft FastTrack Center,", `
                         2)

# --- Change 4: "customers tenant." -> "customer tenant." (Step 4) ---
$d.Content.Find.Execute("customers tenant. There will be Caching servers", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "customer tenant. There will be Caching servers", 2)
